# feat: add 2022-Q1 data
#
# 1. Duplicate the "2021-Q4" sheet (same column layout) to create a new
#    "2022-Q1" sheet, positioned right before "总计", and update its fund
#    data row.
# 2. Insert a new top row into "总计" summarizing the 2022-Q1 quarter,
#    pushing the existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet from the "2021-Q4" template
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($wb.Worksheets.Item("总计"))

$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# Helper: write a value into a cell as TEXT (matches the source data,
# which stores these numeric-looking figures as strings) without
# attaching a new/different cell style to the destination.
$scratch = $newSheet.Cells.Item(60, 60)
$scratch.NumberFormat = "@"

function Write-TextValue($cell, $text) {
    $scratch.Value = $text
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
}

Write-TextValue $newSheet.Cells.Item(2, 4) "0.29"
Write-TextValue $newSheet.Cells.Item(2, 5) "93.32"
Write-TextValue $newSheet.Cells.Item(2, 6) "8.50"
Write-TextValue $newSheet.Cells.Item(2, 7) "0.0246"
$newSheet.Cells.Item(2, 8).Value = 4

$scratch.Delete()

# ---------------------------------------------------------------------
# Step 2: insert the 2022-Q1 summary row at the top of "总计"
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 1
$total.Cells.Item(2, 4).Value = 0.02

# Re-apply the label-column style (bold/bordered/centered) that A3:A7
# already carry, matching the rest of the column.
$total.Cells.Item(3, 1).Copy()
$total.Cells.Item(2, 1).PasteSpecial(-4122)
$total.Application.CutCopyMode = $false

# The inserted row shifted the old rows down but kept their old running
# index (column A); renumber rows 3..7 back into a contiguous 1..5
# sequence so it keeps counting on from the new row 2 (=0).
for ($r = 3; $r -le 7; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}
